$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.889.16"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "2.107.44"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.18"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.94"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.391"
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0782"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  +3.70%  "
$ws.Range("D12").Value = "2.406.56"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.32"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.778"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "2.106.83"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "37.879.56"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.25"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.48"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.05"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.140"
$ws.Range("E27").Value = "  +10.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.00"
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.57"
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.68"
$ws.Range("E32").Value = "  +5.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.45"
$ws.Range("E36").Value = "  +5.63%  "
$ws.Range("E37").Value = "  +4.53%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("E40").Value = "  +6.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("E43").Value = "  +2.39%  "
$ws.Range("D44").Value = "1.459.48"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  +4.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.79"
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.09"
$ws.Range("E48").Value = "  -7.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.37"
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").Value = "2.301.87"
$ws.Range("E51").Value = "  +2.39%  "
